$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Correct the explanation text for the "male" metadata row: was "0 = yes, 1 = no"
$ws.Range("C10").Value = "0 = no, 1 = yes"

# Update the remembered selection on the sheet view
$ws.Range("C23").Select()
